$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, D (Fecha serial), M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @(
    ,@(2, 44874, 240, 29000, 30000, 29500, 1475)
    ,@(3, 44809, 60, 27000, 28000, 27500, 1375)
    ,@(4, 44442, 140, 20000, 21000, 20500, 1025)
    ,@(5, 44462, 100, 19500, 20000, 19750, 988)
    ,@(6, 44434, 100, 20000, 21000, 20500, 1025)
    ,@(7, 44448, 100, 20000, 21000, 20500, 1025)
    ,@(8, 44782, 200, 23500, 24000, 23750, 1188)
    ,@(9, 44335, 200, 19000, 20000, 19500, 975)
    ,@(10, 44428, 100, 20000, 21000, 20500, 1025)
    ,@(11, 44445, 160, 20000, 21000, 20500, 1025)
    ,@(12, 44336, 100, 19500, 20000, 19750, 988)
    ,@(13, 44431, 160, 21000, 22000, 21500, 1075)
    ,@(14, 44879, 100, 28000, 30000, 29000, 1450)
    ,@(15, 44880, 100, 28000, 30000, 29000, 1450)
    ,@(16, 44417, 160, 20000, 21000, 20500, 1025)
    ,@(17, 44810, 100, 27000, 28000, 27500, 1375)
    ,@(18, 44407, 160, 20000, 21000, 20500, 1025)
    ,@(19, 44301, 100, 18000, 19000, 18500, 925)
    ,@(20, 44441, 160, 20000, 21000, 20500, 1025)
    ,@(21, 44420, 160, 20000, 21000, 20500, 1025)
    ,@(22, 44343, 100, 19500, 20000, 19750, 988)
    ,@(23, 44365, 100, 20000, 21000, 20500, 1025)
    ,@(24, 44410, 200, 20000, 21000, 20500, 1025)
    ,@(25, 44882, 120, 28000, 30000, 29000, 1450)
    ,@(26, 44350, 160, 19000, 20000, 19500, 975)
    ,@(27, 44326, 160, 19500, 20000, 19750, 988)
    ,@(28, 44435, 260, 20000, 22000, 21115, 1056)
    ,@(29, 44784, 160, 27000, 28000, 27500, 1375)
    ,@(30, 44473, 40, 19500, 20000, 19750, 988)
    ,@(31, 44418, 200, 20000, 21000, 20500, 1025)
    ,@(32, 44778, 100, 23000, 24000, 23500, 1175)
    ,@(33, 44474, 200, 19000, 20000, 19500, 975)
    ,@(34, 44333, 100, 19500, 20000, 19750, 988)
    ,@(35, 44427, 200, 20000, 21000, 20500, 1025)
    ,@(36, 44776, 160, 23000, 24000, 23500, 1175)
    ,@(37, 44315, 100, 20000, 21000, 20500, 1025)
    ,@(38, 44466, 100, 20000, 21000, 20500, 1025)
    ,@(39, 44781, 160, 23000, 24000, 23500, 1175)
    ,@(40, 44467, 200, 20000, 21000, 20500, 1025)
    ,@(41, 44364, 140, 20000, 21000, 20500, 1025)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D Fecha
    $ws.Cells.Item($r, 13).Value = $row[2]  # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[3]  # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[4]  # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[5]  # P Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[6]  # S Precio $/Kg
}
